$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 14999.667
$ws.Range("J69").Value = 14999.667
$ws.Range("L69").Value = 44999.001
$ws.Range("N69").Value = -46747.001
$ws.Range("H72").Value = 14999.667
$ws.Range("J72").Value = 14999.667
$ws.Range("L72").Value = 134997.003
$ws.Range("N72").Value = -143733.003
$ws.Range("H112").Value = 4164.25
$ws.Range("J112").Value = 4353.8076
$ws.Range("L112").Value = 13061.4228
$ws.Range("N112").Value = -15277.4228
$ws.Range("H113").Value = 4333.3335
$ws.Range("J113").Value = 4333.3335
$ws.Range("L113").Value = 4333.3335
$ws.Range("N113").Value = -10841.3335
$ws.Range("H136").Value = 137872
$ws.Range("J136").Value = 137872
$ws.Range("L136").Value = 137872
$ws.Range("N136").Value = -148072
$ws.Range("H138").Value = 4644.6777
$ws.Range("I138").Value = 1124.3125
$ws.Range("J138").Value = 5954.5815
$ws.Range("K138").Value = 3372.9375
$ws.Range("L138").Value = 17863.7445
$ws.Range("M138").Value = 1767.0625
$ws.Range("N138").Value = -28143.7445
$ws.Range("H140").Value = 59683.8
$ws.Range("J140").Value = 58458.777
$ws.Range("L140").Value = 58458.777
$ws.Range("N140").Value = -68818.777

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6632.8335
$ws.Range("I32").Value = 3701.3635
$ws.Range("K32").Value = 3701.3635
$ws.Range("M32").Value = -3414.3635
$ws.Range("H61").Value = 14019.538
$ws.Range("I61").Value = 27010.5
$ws.Range("J61").Value = 2884.4285
$ws.Range("K61").Value = 27010.5
$ws.Range("L61").Value = 2884.4285
$ws.Range("M61").Value = -26798.5
$ws.Range("N61").Value = -3308.4285
$ws.Range("H74").Value = 9616982
$ws.Range("I74").Value = 11905940
$ws.Range("K74").Value = 11905940
$ws.Range("M74").Value = -11905066
$ws.Range("H77").Value = 9616982
$ws.Range("I77").Value = 11905940
$ws.Range("K77").Value = 59529700
$ws.Range("M77").Value = -59525332
$ws.Range("H80").Value = 69650
$ws.Range("J80").Value = 69800
$ws.Range("L80").Value = 69800
$ws.Range("N80").Value = -71796
$ws.Range("H83").Value = 69650
$ws.Range("J83").Value = 69800
$ws.Range("L83").Value = 209400
$ws.Range("N83").Value = -219384
$ws.Range("H97").Value = 974.06665
$ws.Range("I97").Value = 1201
$ws.Range("K97").Value = 1201
$ws.Range("M97").Value = -705
$ws.Range("H132").Value = 5406.3438
$ws.Range("I132").Value = 2086.5
$ws.Range("J132").Value = 12710
$ws.Range("K132").Value = 6259.5
$ws.Range("L132").Value = 38130
$ws.Range("M132").Value = -3729.5
$ws.Range("N132").Value = -43190
$ws.Range("H136").Value = 14019.538
$ws.Range("I136").Value = 27010.5
$ws.Range("J136").Value = 2884.4285
$ws.Range("K136").Value = 81031.5
$ws.Range("L136").Value = 8653.2855
$ws.Range("M136").Value = -78481.5
$ws.Range("N136").Value = -13753.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4179.8184
$ws.Range("I134").Value = 1997.25
$ws.Range("K134").Value = 5991.75
$ws.Range("M134").Value = -3456.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5905.3477
$ws.Range("I31").Value = 1939.35
$ws.Range("K31").Value = 1939.35
$ws.Range("M31").Value = -1644.35
$ws.Range("H34").Value = 5905.3477
$ws.Range("I34").Value = 1939.35
$ws.Range("K34").Value = 1939.35
$ws.Range("M34").Value = -1737.35
$ws.Range("H94").Value = 2729.15
$ws.Range("J94").Value = 2715.7273
$ws.Range("L94").Value = 2715.7273
$ws.Range("N94").Value = -3617.7273
$ws.Range("H99").Value = 6970.923
$ws.Range("I99").Value = 2770.3333
$ws.Range("K99").Value = 2770.3333
$ws.Range("M99").Value = -1272.3333
$ws.Range("H126").Value = 6970.923
$ws.Range("I126").Value = 2770.3333
$ws.Range("K126").Value = 8310.999899999999
$ws.Range("M126").Value = -5840.999899999999
$ws.Range("H141").Value = 96609.75
$ws.Range("J141").Value = 96609.75
$ws.Range("L141").Value = 96609.75
$ws.Range("N141").Value = -106969.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 47871110
$ws.Range("I4").Value = 1823826.8
$ws.Range("J4").Value = 178338420
$ws.Range("K4").Value = 5471480.4
$ws.Range("L4").Value = 535015260
$ws.Range("M4").Value = -5471368.4
$ws.Range("N4").Value = -535015484
$ws.Range("H33").Value = 354.81818
$ws.Range("J33").Value = 368
$ws.Range("L33").Value = 2208
$ws.Range("N33").Value = -2774
$ws.Range("H61").Value = 485.45456
$ws.Range("I61").Value = 249.2
$ws.Range("J61").Value = 682.3333
$ws.Range("K61").Value = 747.5999999999999
$ws.Range("L61").Value = 2046.9999
$ws.Range("M61").Value = -532.5999999999999
$ws.Range("N61").Value = -2476.9999
$ws.Range("H63").Value = 19000
$ws.Range("J63").Value = 25000
$ws.Range("L63").Value = 75000
$ws.Range("N63").Value = -76498
$ws.Range("H66").Value = 19000
$ws.Range("J66").Value = 25000
$ws.Range("L66").Value = 225000
$ws.Range("N66").Value = -232488
$ws.Range("H81").Value = 4086.2222
$ws.Range("I81").Value = 2569.5
$ws.Range("K81").Value = 7708.5
$ws.Range("M81").Value = -6585.5
$ws.Range("H84").Value = 4086.2222
$ws.Range("I84").Value = 2569.5
$ws.Range("K84").Value = 23125.5
$ws.Range("M84").Value = -17509.5
$ws.Range("H87").Value = 20000
$ws.Range("I87").Value = 10000
$ws.Range("J87").Value = 25000
$ws.Range("K87").Value = 30000
$ws.Range("L87").Value = 75000
$ws.Range("M87").Value = -28752
$ws.Range("N87").Value = -77496
$ws.Range("H90").Value = 20000
$ws.Range("I90").Value = 10000
$ws.Range("J90").Value = 25000
$ws.Range("K90").Value = 90000
$ws.Range("L90").Value = 225000
$ws.Range("M90").Value = -83760
$ws.Range("N90").Value = -237480

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 769763.3
$ws.Range("I80").Value = 2082724.8
$ws.Range("K80").Value = 2082724.8
$ws.Range("M80").Value = -2081726.8
$ws.Range("H83").Value = 769763.3
$ws.Range("I83").Value = 2082724.8
$ws.Range("K83").Value = 10413624
$ws.Range("M83").Value = -10408632
$ws.Range("H126").Value = 4242.852
$ws.Range("I126").Value = 3102.2856
$ws.Range("J126").Value = 5471.154
$ws.Range("K126").Value = 9306.856800000001
$ws.Range("L126").Value = 16413.462
$ws.Range("M126").Value = -6836.856800000001
$ws.Range("N126").Value = -21353.462
$ws.Range("H132").Value = 6706
$ws.Range("I132").Value = 7877.6
$ws.Range("J132").Value = 5534.4
$ws.Range("K132").Value = 23632.8
$ws.Range("L132").Value = 16603.2
$ws.Range("M132").Value = -21102.8
$ws.Range("N132").Value = -21663.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7618.095
$ws.Range("J46").Value = 7618.095
$ws.Range("L46").Value = 7618.095
$ws.Range("N46").Value = -7994.095
$ws.Range("H82").Value = 15625750
$ws.Range("I82").Value = 31250000
$ws.Range("J82").Value = 1500
$ws.Range("K82").Value = 31250000
$ws.Range("L82").Value = 1500
$ws.Range("M82").Value = -31249639
$ws.Range("N82").Value = -2222
$ws.Range("H85").Value = 15625750
$ws.Range("I85").Value = 31250000
$ws.Range("J85").Value = 1500
$ws.Range("K85").Value = 31250000
$ws.Range("L85").Value = 1500
$ws.Range("M85").Value = -31248752
$ws.Range("N85").Value = -3996
$ws.Range("H136").Value = 2129.425
$ws.Range("I136").Value = 1956.0857
$ws.Range("J136").Value = 3342.8
$ws.Range("K136").Value = 5868.257100000001
$ws.Range("L136").Value = 10028.4
$ws.Range("M136").Value = -3318.257100000001
$ws.Range("N136").Value = -15128.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 65000
$ws.Range("J49").Value = 65000
$ws.Range("L49").Value = 65000
$ws.Range("N49").Value = -65460
$ws.Range("H81").Value = 10423166
$ws.Range("I81").Value = 10422666
$ws.Range("K81").Value = 20845332
$ws.Range("M81").Value = -20844271
$ws.Range("H84").Value = 10423166
$ws.Range("I84").Value = 10422666
$ws.Range("K84").Value = 104226660
$ws.Range("M84").Value = -104221356
